$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same event rows and need
# their "想去人数" (F column) counts bumped.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 343
    $ws.Range("F4").Value = 73
    $ws.Range("F5").Value = 291
}
